$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.002.22'
$ws.Range('E2').Value = '  +1.49%  '
$ws.Range('D3').Value = '2.912.88'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '199.80'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +5.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '598.54'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.550'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.197'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('D10').Value = '2.911.67'
$ws.Range('E10').Value = '  +2.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.427'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +14.75%  '
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.89'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('D14').Value = '3.448.37'
$ws.Range('E14').Value = '  +2.31%  '
$ws.Range('D15').Value = '75.845.27'
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.42'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').Value = '2.911.63'
$ws.Range('E18').Value = '  +2.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.90'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('E20').Value = '  +2.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '377.77'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('E23').Value = '  +1.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.38'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('D26').Value = '3.062.92'
$ws.Range('E26').Value = '  +2.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.21'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.73'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000110'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +5.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '505.03'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.95%  '
$ws.Range('E33').Value = '  -2.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '164.26'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.19'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.70'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.07%  '
$ws.Range('E39').Value = '  -6.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '180.08'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.25%  '
$ws.Range('E42').Value = '  +0.38%  '
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0924'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +7.85%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.66'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.18'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.54%  '
$ws.Range('E47').Value = '  -3.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.35'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.661'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +7.08%  '
$ws.Range('E51').Value = '  -0.89%  '
